# Music_Zone.pptx - slide with the project goal statement.
#
# Original single run:
#   "разработать удобный в использовании каталог музыкального приложения для поиска и прослушивания музыки "
#
# becomes three runs, with "удобный в использовании " removed entirely:
#   "разработать " + "каталог " + "музыкального приложения для поиска и прослушивания музыки "

$p = $ppt.ActivePresentation

$oldPhrase = "разработать удобный в использовании каталог музыкального приложения для поиска и прослушивания музыки "
$keepStart = "разработать "
$removedPhrase = "удобный в использовании "
$keepMiddle = "каталог "

# Locate the shape/paragraph that holds the sentence, rather than assuming a
# fixed slide/shape number.
$targetShape = $null
$targetSlide = $null
$targetParaIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.Contains($oldPhrase)) {
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text.Contains($oldPhrase)) {
                        $targetSlide = $slide
                        $targetShape = $shape
                        $targetParaIndex = $pi
                    }
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $para = $tr.Paragraphs($targetParaIndex, 1)
    $paraText = $para.Text

    $fullStart = $paraText.IndexOf($oldPhrase)
    # 1-indexed character position where the target phrase begins.
    $baseStart = $fullStart + 1

    # Drop "удобный в использовании " (it sits right after "разработать ").
    $removedStart = $baseStart + $keepStart.Length
    $removedRange = $para.Characters($removedStart, $removedPhrase.Length)
    $removedRange.Text = ""

    # "каталог " now begins immediately after "разработать " and becomes its
    # own run, splitting the remaining tail ("музыкального приложения...")
    # into a third run.
    $midStart = $baseStart + $keepStart.Length
    $midRange = $para.Characters($midStart, $keepMiddle.Length)
    $midRange.Text = $keepMiddle
}
